# netCrypto.xlsx update — "Add files via upload"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount figure in T2
$ws.Range("T2").Value = 70297

# Scroll the view so column L is the left-most visible column (mirrors
# topLeftCell moving from H1 to L1 in the saved view state)
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1

# Move the active selection to R13 (was T3)
$ws.Range("R13").Select()
